# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
#
# Updates the "Metadata" sheet of the StructureDefinition-process-record-id
# workbook:
#   - Version bumped 5.0.0 -> 6.0.0
#   - Date bumped to the new publication date
#   - Publisher value filled in ("Alvearie Team")
#   - The row that used to read "Contact" / "No display for ContactDetail"
#     now reads "Jurisdiction" / "United States of America"
#   - The second, now-redundant "Contact" row is removed entirely, which
#     shifts every following row up by one (table shrinks from 21 to 20 rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: row 3
$ws.Range("B3").Value = "6.0.0"

# Date: row 8
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher: row 9 (label stays the same, value gets filled in)
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 used to be the first "Contact" row; repurpose it as "Jurisdiction"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was the duplicate "Contact" / "No display for ContactDetail" row -
# delete it outright so every subsequent row (Description, Purpose, ...,
# Context) shifts up by one.
$ws.Rows.Item(11).Delete()
